# The deck's slide-master theme (ppt/theme/theme1.xml, "Integral") is being
# swapped for the default Office theme's colour scheme ("Office Theme" /
# "Office": dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).
#
# PowerPoint's ThemeColorScheme object exposes exactly those 12 DrawingML
# theme colours (in clrScheme document order) through any slide, so we
# rewrite them in place via the COM object model - no direct file/XML
# access is needed (or available) in this host.

function HexToRGB($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Index order matches <a:clrScheme>: dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink.
$tcs.Colors(1).RGB  = HexToRGB "000000"   # dk1
$tcs.Colors(2).RGB  = HexToRGB "FFFFFF"   # lt1
$tcs.Colors(3).RGB  = HexToRGB "44546A"   # dk2
$tcs.Colors(4).RGB  = HexToRGB "E7E6E6"   # lt2
$tcs.Colors(5).RGB  = HexToRGB "5B9BD5"   # accent1
$tcs.Colors(6).RGB  = HexToRGB "ED7D31"   # accent2
$tcs.Colors(7).RGB  = HexToRGB "A5A5A5"   # accent3
$tcs.Colors(8).RGB  = HexToRGB "FFC000"   # accent4
$tcs.Colors(9).RGB  = HexToRGB "4472C4"   # accent5
$tcs.Colors(10).RGB = HexToRGB "70AD47"   # accent6
$tcs.Colors(11).RGB = HexToRGB "0563C1"   # hlink
$tcs.Colors(12).RGB = HexToRGB "954F72"   # folHlink
